# Apply the "results" sheet changes:
# 1) Split the combined D1 header into separate D1/E1 headers.
# 2) Update the numeric values for the "Расстояние до переноса стиля arcface" (D)
#    and "Расстояние после переноса стиля arcface" (E) columns (and a few
#    slightly-adjusted C values) to reflect the new embedding calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("results")

# --- Header split ---
$ws.Range("D1").Value = "Расстояние до переноса стиля arcface"
$ws.Range("E1").Value = "Расстояние после переноса стиля arcface"

# --- Updated values ---
# Row 2
$ws.Range("C2").Value = 0.6900752004639339
$ws.Range("D2").Value = 0.008492808789014816
$ws.Range("E2").Value = 0.002082576043903828

# Row 3
$ws.Range("D3").Value = 0.01410352624952793
$ws.Range("E3").Value = 0.01060810126364231

# Row 4
$ws.Range("D4").Value = 0.01806437224149704
$ws.Range("E4").Value = 0.003665033960714936

# Row 5
$ws.Range("C5").Value = 0.7964303098175904
$ws.Range("D5").Value = 0.009923557750880718
$ws.Range("E5").Value = 0.005861029028892517

# Row 6
$ws.Range("D6").Value = 0.01835355162620544
$ws.Range("E6").Value = 0.01003515347838402

# Row 7
$ws.Range("D7").Value = 0.01772580668330193
$ws.Range("E7").Value = 0.01566718332469463

# Row 8
$ws.Range("D8").Value = 0.01674813218414783
$ws.Range("E8").Value = 0.00910513661801815

# Row 9
$ws.Range("D9").Value = 0.02624757029116154
$ws.Range("E9").Value = 0.01961797103285789

# Row 10
$ws.Range("C10").Value = 0.7330400849497321
$ws.Range("D10").Value = 0.0465083047747612
$ws.Range("E10").Value = 0.03962106630206108

# Row 11
$ws.Range("C11").Value = 0.754820729126691
$ws.Range("D11").Value = 0.02004562132060528
$ws.Range("E11").Value = 0.01289359852671623

# Row 12
$ws.Range("C12").Value = 0.4725797102277756
$ws.Range("D12").Value = 0.006884787697345018
$ws.Range("E12").Value = 0.005170291755348444

# Row 13
$ws.Range("D13").Value = 0.0232948437333107
$ws.Range("E13").Value = 0.01254448667168617

# Row 14
$ws.Range("D14").Value = 0.02751919999718666
$ws.Range("E14").Value = 0.02013404667377472

# Row 15
$ws.Range("C15").Value = 0.6791948264882282
$ws.Range("D15").Value = 0.03104481101036072
$ws.Range("E15").Value = 0.01083913818001747

# Row 16
$ws.Range("D16").Value = 0.02553591132164001
$ws.Range("E16").Value = 0.01203899551182985

# Row 17
$ws.Range("D17").Value = 0.01383387669920921
$ws.Range("E17").Value = 0.01148590538650751

# Row 18
$ws.Range("C18").Value = 0.6392436588714024
$ws.Range("D18").Value = 0.1120398044586182
$ws.Range("E18").Value = 0.07538660615682602

# Row 19
$ws.Range("C19").Value = 0.6552471888204152
$ws.Range("D19").Value = 0.02931328490376472
$ws.Range("E19").Value = 0.02313164621591568
